# Week 44 profile updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two profiles flipped from "not yet updated" (0) to "updated" (1)
$ws.Range("C37").Value = 1
$ws.Range("C57").Value = 1

# New profile entry appended on the first free row (115).
# Write the username (column B) before the full name (column A) so the
# shared-string table picks up "archit_desai" ahead of "Archit Anuj Desai",
# matching how the roster has always grown one profile at a time.
$ws.Range("B115").Value = "archit_desai"
$ws.Range("A115").Value = "Archit Anuj Desai"
$ws.Range("C115").Value = 0

# Reuse the same look as the most recent "new profile" rows in this sheet
# (name cell + checkbox cell carry a distinct format, username stays plain).
$ws.Range("A113").Copy()
$ws.Range("A115").PasteSpecial(-4122)
$ws.Range("C109").Copy()
$ws.Range("C115").PasteSpecial(-4122)

# Leave the view where the reviewer was working (around the two toggles).
$ws.Range("C57").Select()
